# Add Spring related documentation
# Locates the three blank paragraphs that follow the "@Autowired" /
# kipalog link block (just before the "Other notes" heading) and:
#   1. turns the first blank paragraph into a bold "Spring documentation"
#      heading
#   2. turns the second blank paragraph into a hyperlink pointing at the
#      Spring Data CrudRepository#existsById documentation
#   3. leaves the third blank paragraph untouched
#   4. inserts two more blank paragraphs before "Other notes"

$d = $word.ActiveDocument

# Locate the "@Autowired" run as an anchor, then walk forward from there.
$found = $d.Content.Find
$found.ClearFormatting()
$anchor = $d.Content
$anchor.Find.Execute("https://kipalog.com/posts/Spring-Boot--1--Hu-o--ng-dan--Component-va---Autowired") | Out-Null

# Paragraph that holds the kipalog URL text we just located.
$urlParaRange = $anchor.Paragraphs.Item(1).Range
$urlParaIndex = $urlParaRange.Paragraphs.Item(1).Range.Start

# Walk the paragraph collection to find the index of that paragraph, then
# use the three following blank paragraphs.
$allParas = $d.Paragraphs
$urlIndex = 0
for ($i = 1; $i -le $allParas.Count; $i++) {
    if ($allParas.Item($i).Range.Start -eq $urlParaRange.Start) {
        $urlIndex = $i
        break
    }
}

$headingIndex = $urlIndex + 1
$linkIndex = $urlIndex + 2
$blankIndex = $urlIndex + 3

# --- 1. "Spring documentation" heading -----------------------------------
$headingPara = $d.Paragraphs.Item($headingIndex)
$headingRange = $headingPara.Range
$headingRange.Text = "Spring documentation"
$headingRange2 = $d.Paragraphs.Item($headingIndex).Range
$headingRange2.Font.Name = "Consolas"
$headingRange2.Font.Size = 9
$headingRange2.Font.Color = 0
$headingRange2.Font.Bold = 1
$headingRange2.Font.BoldBi = 1

# --- 2. Hyperlink to the CrudRepository docs ------------------------------
$linkPara = $d.Paragraphs.Item($linkIndex)
$linkRange = $linkPara.Range
$hyperlink = $d.Hyperlinks.Add(
    $linkRange,
    "https://docs.spring.io/spring-data/commons/docs/current/api/org/springframework/data/repository/CrudRepository.html?is-external=true#existsById-ID-",
    "",
    "",
    "https://docs.spring.io/spring-data/commons/docs/current/api/org/springframework/data/repository/CrudRepository.html?is-external=true#existsById-ID-"
)
$hyperlinkRange = $hyperlink.Range
$hyperlinkRange.Font.Name = "Consolas"
$hyperlinkRange.Font.Size = 9

# --- 3. leave the 3rd blank paragraph (index $blankIndex) untouched ------

# --- 4. insert two more blank paragraphs before "Other notes" ------------
$blankPara = $d.Paragraphs.Item($blankIndex)
$insertPoint = $d.Range($blankPara.Range.Start, $blankPara.Range.Start)
$insertPoint.InsertParagraphAfter()
$blankPara2 = $d.Paragraphs.Item($blankIndex)
$insertPoint2 = $d.Range($blankPara2.Range.Start, $blankPara2.Range.Start)
$insertPoint2.InsertParagraphAfter()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
